{"js": "// Replace each three-digit-by-one-digit division equation's result text\n// with its updated value, per the commit diff. Each old string is unique\n// within the document, so body.search(...) safely targets exactly one run.\nconst replacements = [\n  [\"420\u00f75=84, 0\", \"370\u00f73=123, 1\"],\n  [\"190\u00f76=31, 4\", \"106\u00f74=26, 2\"],\n  [\"320\u00f79=35, 5\", \"750\u00f78=93, 6\"],\n  [\"366\u00f72=183, 0\", \"513\u00f75=102, 3\"],\n  [\"577\u00f75=115, 2\", \"560\u00f78=70, 0\"],\n  [\"450\u00f74=112, 2\", \"499\u00f72=249, 1\"],\n  [\"762\u00f77=108, 6\", \"994\u00f72=497, 0\"],\n  [\"946\u00f77=135, 1\", \"777\u00f76=129, 3\"],\n  [\"387\u00f74=96, 3\", \"809\u00f77=115, 4\"],\n  [\"442\u00f74=110, 2\", \"962\u00f76=160, 2\"],\n  [\"388\u00f73=129, 1\", \"243\u00f76=40, 3\"],\n  [\"723\u00f75=144, 3\", \"631\u00f74=157, 3\"],\n  [\"584\u00f77=83, 3\", \"969\u00f78=121, 1\"],\n  [\"643\u00f74=160, 3\", \"413\u00f72=206, 1\"],\n  [\"183\u00f79=20, 3\", \"260\u00f74=65, 0\"],\n  [\"133\u00f73=44, 1\", \"773\u00f78=96, 5\"],\n  [\"781\u00f72=390, 1\", \"382\u00f72=191, 0\"],\n  [\"981\u00f76=163, 3\", \"283\u00f72=141, 1\"],\n  [\"154\u00f73=51, 1\", \"949\u00f72=474, 1\"],\n  [\"212\u00f73=70, 2\", \"743\u00f79=82, 5\"],\n  [\"525\u00f79=58, 3\", \"743\u00f76=123, 5\"],\n  [\"891\u00f72=445, 1\", \"550\u00f74=137, 2\"],\n  [\"671\u00f76=111, 5\", \"983\u00f73=327, 2\"],\n  [\"118\u00f72=59, 0\", \"753\u00f78=94, 1\"],\n  [\"877\u00f72=438, 1\", \"440\u00f78=55, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update each three-digit-by-one-digit division equation result in the table\n# to its new value, per the commit diff. Each \"old\" string is unique in the\n# document, so a plain Find/Replace (MatchWholeWord off, no wildcards) safely\n# retargets exactly one run per pair without touching formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"420\u00f75=84, 0\", \"370\u00f73=123, 1\"),\n    @(\"190\u00f76=31, 4\", \"106\u00f74=26, 2\"),\n    @(\"320\u00f79=35, 5\", \"750\u00f78=93, 6\"),\n    @(\"366\u00f72=183, 0\", \"513\u00f75=102, 3\"),\n    @(\"577\u00f75=115, 2\", \"560\u00f78=70, 0\"),\n    @(\"450\u00f74=112, 2\", \"499\u00f72=249, 1\"),\n    @(\"762\u00f77=108, 6\", \"994\u00f72=497, 0\"),\n    @(\"946\u00f77=135, 1\", \"777\u00f76=129, 3\"),\n    @(\"387\u00f74=96, 3\", \"809\u00f77=115, 4\"),\n    @(\"442\u00f74=110, 2\", \"962\u00f76=160, 2\"),\n    @(\"388\u00f73=129, 1\", \"243\u00f76=40, 3\"),\n    @(\"723\u00f75=144, 3\", \"631\u00f74=157, 3\"),\n    @(\"584\u00f77=83, 3\", \"969\u00f78=121, 1\"),\n    @(\"643\u00f74=160, 3\", \"413\u00f72=206, 1\"),\n    @(\"183\u00f79=20, 3\", \"260\u00f74=65, 0\"),\n    @(\"133\u00f73=44, 1\", \"773\u00f78=96, 5\"),\n    @(\"781\u00f72=390, 1\", \"382\u00f72=191, 0\"),\n    @(\"981\u00f76=163, 3\", \"283\u00f72=141, 1\"),\n    @(\"154\u00f73=51, 1\", \"949\u00f72=474, 1\"),\n    @(\"212\u00f73=70, 2\", \"743\u00f79=82, 5\"),\n    @(\"525\u00f79=58, 3\", \"743\u00f76=123, 5\"),\n    @(\"891\u00f72=445, 1\", \"550\u00f74=137, 2\"),\n    @(\"671\u00f76=111, 5\", \"983\u00f73=327, 2\"),\n    @(\"118\u00f72=59, 0\", \"753\u00f78=94, 1\"),\n    @(\"877\u00f72=438, 1\", \"440\u00f78=55, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n\n    $found = $rng.Find.Execute(\n        $oldText,  # FindText\n        $false,    # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $newText,  # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
